$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('Sending cluster', 'Ligand symbol', 'Receptor symbol', 'Target cluster', 'Ligand-expressing cells', 'Ligand detection rate', 'Ligand average expression value', 'Ligand total expression value', 'Ligand derived specificity of average expression value', 'Ligand derived specificity of total expression value', 'Receptor-expressing cells', 'Receptor detection rate', 'Receptor average expression value', 'Receptor total expression value', 'Receptor derived specificity of average expression value', 'Receptor derived specificity of total expression value', 'Edge average expression weight', 'Edge total expression weight', 'Edge average expression derived specificity', 'Edge total expression derived specificity'),
  @('ECs', 'Cd38', 'Pecam1', 'ECs', 3, 1, 38.12257066666666, 114.367712, 0.696502358511838, 0.6965023585118381, 3, 1, 218.9722396666667, 656.9167190000001, 0.9460661569256378, 0.9460661569256379, 8347.784680730769, 75130.06212657693, 0.6589373096069373, 0.6589373096069375),
  @('ECs', 'Cd38', 'Pecam1', 'FAPs', 3, 1, 38.12257066666666, 114.367712, 0.696502358511838, 0.6965023585118381, 3, 1, 5.232695, 15.698085, 0.02260777739018391, 0.02260777739018392, 199.4837849146133, 1795.35406423152, 0.0157463702729737, 0.01574637027297371),
  @('ECs', 'Cd38', 'Pecam1', 'M2', 3, 1, 38.12257066666666, 114.367712, 0.696502358511838, 0.6965023585118381, 3, 1, 3.718869666666667, 11.156609, 0.01606731857429249, 0.0160673185742925, 141.7728716676231, 1275.955845008608, 0.01119092528195578, 0.01119092528195579),
  @('ECs', 'Cd38', 'Pecam1', 'sCs', 3, 1, 38.12257066666666, 114.367712, 0.696502358511838, 0.6965023585118381, 3, 1, 3.531721333333333, 10.595164, 0.01525874710988573, 0.01525874710988573, 134.6382961049742, 1211.744664944768, 0.0106277533499711, 0.01062775334997111),
  @('FAPs', 'Cd38', 'Pecam1', 'ECs', 3, 1, 2.823821666666666, 8.471464999999998, 0.05159144350593012, 0.05159144350593013, 3, 1, 218.9722396666667, 656.9167190000001, 0.9460661569256378, 0.9460661569256379, 618.3385547692594, 5565.046992923334, 0.04880891868790146, 0.04880891868790148),
  @('FAPs', 'Cd38', 'Pecam1', 'FAPs', 3, 1, 2.823821666666666, 8.471464999999998, 0.05159144350593012, 0.05159144350593013, 3, 1, 5.232695, 15.698085, 0.02260777739018391, 0.02260777739018392, 14.77619751605833, 132.985777644525, 0.001166367870020318, 0.001166367870020318),
  @('FAPs', 'Cd38', 'Pecam1', 'M2', 3, 1, 2.823821666666666, 8.471464999999998, 0.05159144350593012, 0.05159144350593013, 3, 1, 3.718869666666667, 11.156609, 0.01606731857429249, 0.0160673185742925, 10.50142474024278, 94.51282266218499, 0.0008289361585173928, 0.0008289361585173932),
  @('FAPs', 'Cd38', 'Pecam1', 'sCs', 3, 1, 2.823821666666666, 8.471464999999998, 0.05159144350593012, 0.05159144350593013, 3, 1, 3.531721333333333, 10.595164, 0.01525874710988573, 0.01525874710988573, 9.972951221695554, 89.75656099525999, 0.0007872207894909443, 0.0007872207894909446),
  @('M2', 'Cd38', 'Pecam1', 'ECs', 3, 1, 13.44991933333333, 40.349758, 0.2457310819716487, 0.2457310819716487, 3, 1, 218.9722396666667, 656.9167190000001, 0.9460661569256378, 0.9460661569256379, 2945.158959756001, 26506.43063780401, 0.2324778603580965, 0.2324778603580966),
  @('M2', 'Cd38', 'Pecam1', 'FAPs', 3, 1, 13.44991933333333, 40.349758, 0.2457310819716487, 0.2457310819716487, 3, 1, 5.232695, 15.698085, 0.02260777739018391, 0.02260777739018392, 70.37932564593666, 633.4139308134299, 0.005555433599064069, 0.005555433599064069),
  @('M2', 'Cd38', 'Pecam1', 'M2', 3, 1, 13.44991933333333, 40.349758, 0.2457310819716487, 0.2457310819716487, 3, 1, 3.718869666666667, 11.156609, 0.01606731857429249, 0.0160673185742925, 50.01849702784689, 450.1664732506221, 0.003948239577644062, 0.003948239577644063),
  @('M2', 'Cd38', 'Pecam1', 'sCs', 3, 1, 13.44991933333333, 40.349758, 0.2457310819716487, 0.2457310819716487, 3, 1, 3.531721333333333, 10.595164, 0.01525874710988573, 0.01525874710988573, 47.50136704114578, 427.5123033703121, 0.003749548436843988, 0.003749548436843989),
  @('sCs', 'Cd38', 'Pecam1', 'ECs', 1, 0.3333333333333333, 0.3379906666666667, 1.013972, 0.006175116010583175, 0.006175116010583175, 3, 1, 218.9722396666667, 656.9167190000001, 0.9460661569256378, 0.9460661569256379, 74.0105732664298, 666.0951593978681, 0.0058420682727024, 0.005842068272702401),
  @('sCs', 'Cd38', 'Pecam1', 'FAPs', 1, 0.3333333333333333, 0.3379906666666667, 1.013972, 0.006175116010583175, 0.006175116010583175, 3, 1, 5.232695, 15.698085, 0.02260777739018391, 0.02260777739018392, 1.768602071513333, 15.91741864362, 0.000139605648125825, 0.000139605648125825),
  @('sCs', 'Cd38', 'Pecam1', 'M2', 1, 0.3333333333333333, 0.3379906666666667, 1.013972, 0.006175116010583175, 0.006175116010583175, 3, 1, 3.718869666666667, 11.156609, 0.01606731857429249, 0.0160673185742925, 1.256943237883111, 11.312489140948, 0.00009921755617525401, 0.00009921755617525403),
  @('sCs', 'Cd38', 'Pecam1', 'sCs', 1, 0.3333333333333333, 0.3379906666666667, 1.013972, 0.006175116010583175, 0.006175116010583175, 3, 1, 3.531721333333333, 10.595164, 0.01525874710988573, 0.01525874710988573, 1.193688847934222, 10.743199631408, 0.00009422453357969514, 0.00009422453357969516)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}
